# Applies the cryptos.xlsx price/volume refresh described in the commit:
# "Updated cryptos list on Sat Nov 16 19:15:21 UTC 2024 with GitHub Actions"
#
# Only the cells that actually changed are touched. The "Price" column (D)
# stores plain-looking numbers (e.g. "216.71") as *text* in the original
# workbook, so a leading apostrophe is used where needed to stop Excel
# from auto-converting the assigned string into a real number (which would
# silently drop meaningful trailing zeros, e.g. "35.20" -> 35.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.224.88"
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").Value = "3.181.00"
$ws.Range("E3").Value = "  +4.98%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'216.71"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").Value = "'628.79"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("E7").Value = "  +31.23%  "

$ws.Range("D8").Value = "'0.372"
$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "3.181.46"
$ws.Range("E10").Value = "  +5.12%  "

$ws.Range("D11").Value = "'0.762"
$ws.Range("E11").Value = "  +14.54%  "

$ws.Range("D12").Value = "'0.202"
$ws.Range("E12").Value = "  +7.43%  "

$ws.Range("E13").Value = "  +2.32%  "

$ws.Range("D14").Value = "'5.67"
$ws.Range("E14").Value = "  +6.81%  "

$ws.Range("D15").Value = "'35.20"
$ws.Range("E15").Value = "  +9.19%  "

$ws.Range("D16").Value = "90.888.26"
$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("D17").Value = "3.760.49"
$ws.Range("E17").Value = "  +5.07%  "

$ws.Range("D18").Value = "3.189.49"
$ws.Range("E18").Value = "  +3.93%  "

$ws.Range("E19").Value = "  +13.20%  "

$ws.Range("D20").Value = "'14.74"
$ws.Range("E20").Value = "  +10.13%  "

$ws.Range("D21").Value = "'476.83"
$ws.Range("E21").Value = "  +12.53%  "

$ws.Range("D22").Value = "'0.0000213"
$ws.Range("E22").Value = "  -2.38%  "

$ws.Range("D23").Value = "'9.23"
$ws.Range("E23").Value = "  +12.06%  "

$ws.Range("D24").Value = "'5.34"
$ws.Range("E24").Value = "  +6.10%  "

$ws.Range("D25").Value = "'96.81"
$ws.Range("E25").Value = "  +17.41%  "

$ws.Range("D26").Value = "'5.77"
$ws.Range("E26").Value = "  +7.70%  "

$ws.Range("D27").Value = "'12.44"
$ws.Range("E27").Value = "  +7.80%  "

$ws.Range("D28").Value = "3.342.22"
$ws.Range("E28").Value = "  +4.84%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  +12.07%  "

$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").Value = "'28.50"
$ws.Range("E33").Value = "  +25.23%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.194"
$ws.Range("E34").Value = "  +42.45%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'529.18"
$ws.Range("E35").Value = "  +5.81%  "

$ws.Range("D36").Value = "'0.146"
$ws.Range("E36").Value = "  +10.17%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.66"
$ws.Range("E37").Value = "  -2.97%  "

$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.94"
$ws.Range("E38").Value = "  +7.72%  "

$ws.Range("D39").Value = "'7.01"
$ws.Range("E39").Value = "  +5.82%  "

$ws.Range("D40").Value = "'1.31"
$ws.Range("E40").Value = "  +5.85%  "

$ws.Range("D41").Value = "'0.0898"
$ws.Range("E41").Value = "  +29.23%  "

$ws.Range("D42").Value = "'22.26"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "'0.422"
$ws.Range("E43").Value = "  +18.78%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  +9.97%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").Value = "'0.713"
$ws.Range("E47").Value = "  +20.85%  "

$ws.Range("D48").Value = "'151.46"
$ws.Range("E48").Value = "  +5.47%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "'1.37"
$ws.Range("E49").Value = "  +12.97%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.59"
$ws.Range("E50").Value = "  +9.45%  "

$ws.Range("D51").Value = "'45.57"
$ws.Range("E51").Value = "  +4.77%  "
